# Bgn-Fgfr3.xlsx -- refresh NATMI ligand-receptor pair table with updated TPM-based
# values (commit: "update scripts wuth new tpm"). The recompute changed per-cluster
# expression stats and, because more sending/target cluster pairs now clear the
# pipelines significance threshold, grew the table from 15 data rows (4 target
# clusters x ~4 senders) to 20 data rows (4 target clusters x 5 senders).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one data row in column order A..T:
#   Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
#   Ligand-expressing cells, Ligand detection rate,
#   Ligand average/total expression value,
#   Ligand derived specificity of average/total expression value,
#   Receptor-expressing cells, Receptor detection rate,
#   Receptor average/total expression value,
#   Receptor derived specificity of average/total expression value,
#   Edge average/total expression weight,
#   Edge average/total expression derived specificity.
$rows = @(
  @("ECs", "Bgn", "Fgfr3", "ECs", 3.0, 1.0, 64.51746766666666, 193.552403, 0.01908790637288361, 0.0192326001326142, 3.0, 1.0, 2.330840333333333, 6.992521, 0.6715345129768794, 0.7003397275969581, 150.3799158419959, 1353.419242577963, 0.01281818790986267, 0.01346935393785625),
  @("ECs", "Bgn", "Fgfr3", "FAPs", 3.0, 1.0, 64.51746766666666, 193.552403, 0.01908790637288361, 0.0192326001326142, 3.0, 1.0, 0.657666, 1.972998, 0.1894790521235985, 0.1976066831789769, 42.43094489046599, 381.8785040141939, 0.003616758406557982, 0.003800490321113444),
  @("ECs", "Bgn", "Fgfr3", "MuSCs", 3.0, 1.0, 64.51746766666666, 193.552403, 0.01908790637288361, 0.0192326001326142, 2.0, 1.0, 0.42828, 0.85656, 0.123391035029171, 0.08578923067523865, 27.63154105228, 165.78924631368, 0.002355276523890018, 0.001649949969261465),
  @("ECs", "Bgn", "Fgfr3", "Resolving-Mac", 3.0, 1.0, 64.51746766666666, 193.552403, 0.01908790637288361, 0.0192326001326142, 1.0, 0.3333333333333333, 0.05413033333333334, 0.162391, 0.01559539987035126, 0.01626435854882633, 3.492352030619222, 31.431168275573, 0.000297683532572946, 0.0003128059043830421),
  @("FAPs", "Bgn", "Fgfr3", "ECs", 3.0, 1.0, 3231.170084666667, 9693.510254, 0.9559623816860557, 0.9632089486203771, 3.0, 1.0, 2.330840333333333, 6.992521, 0.6715345129768794, 0.7003397275969581, 7531.341557201148, 67782.07401481034, 0.6419617324097632, 0.6745734926957473),
  @("FAPs", "Bgn", "Fgfr3", "FAPs", 3.0, 1.0, 3231.170084666667, 9693.510254, 0.9559623816860557, 0.9632089486203771, 3.0, 1.0, 0.657666, 1.972998, 0.1894790521235985, 0.1976066831789769, 2125.030704902388, 19125.27634412149, 0.1811348459476915, 0.1903365255451824),
  @("FAPs", "Bgn", "Fgfr3", "MuSCs", 3.0, 1.0, 3231.170084666667, 9693.510254, 0.9559623816860557, 0.9632089486203771, 2.0, 1.0, 0.42828, 0.85656, 0.123391035029171, 0.08578923067523865, 1383.84552386104, 8303.07314316624, 0.1179571877251938, 0.08263295468164762),
  @("FAPs", "Bgn", "Fgfr3", "Resolving-Mac", 3.0, 1.0, 3231.170084666667, 9693.510254, 0.9559623816860557, 0.9632089486203771, 1.0, 0.3333333333333333, 0.05413033333333334, 0.162391, 0.01559539987035126, 0.01626435854882633, 174.9043137397016, 1574.138823657314, 0.01490861560340739, 0.01566597569779985),
  @("Inflammatory-Mac", "Bgn", "Fgfr3", "ECs", 3.0, 1.0, 4.815307333333333, 14.445922, 0.001424639541189161, 0.001435438863411757, 3.0, 1.0, 2.330840333333333, 6.992521, 0.6715345129768794, 0.7003397275969581, 11.22371254992911, 101.013412949362, 0.0009566946204600682, 0.001005294862583877),
  @("Inflammatory-Mac", "Bgn", "Fgfr3", "FAPs", 3.0, 1.0, 4.815307333333333, 14.445922, 0.001424639541189161, 0.001435438863411757, 3.0, 1.0, 0.657666, 1.972998, 0.1894790521235985, 0.1976066831789769, 3.166863912684, 28.501775214156, 0.0002699393498823205, 0.0002836523127049979),
  @("Inflammatory-Mac", "Bgn", "Fgfr3", "MuSCs", 3.0, 1.0, 4.815307333333333, 14.445922, 0.001424639541189161, 0.001435438863411757, 2.0, 1.0, 0.42828, 0.85656, 0.123391035029171, 0.08578923067523865, 2.06229982472, 12.37379894832, 0.0001757877475308139, 0.0001231451957734336),
  @("Inflammatory-Mac", "Bgn", "Fgfr3", "Resolving-Mac", 3.0, 1.0, 4.815307333333333, 14.445922, 0.001424639541189161, 0.001435438863411757, 1.0, 0.3333333333333333, 0.05413033333333334, 0.162391, 0.01559539987035126, 0.01626435854882633, 0.2606541910557778, 2.345887719502, 0.0000222178233159587109518, 0.0000233464923494485607929),
  @("MuSCs", "Bgn", "Fgfr3", "ECs", 2.0, 1.0, 76.28726950000001, 152.574539, 0.02257007769093324, 0.01516077844306046, 3.0, 1.0, 2.330840333333333, 6.992521, 0.6715345129768794, 0.7003397275969581, 177.8134446704698, 1066.880668022819, 0.01515658613003119, 0.0106176954449708),
  @("MuSCs", "Bgn", "Fgfr3", "FAPs", 2.0, 1.0, 76.28726950000001, 152.574539, 0.02257007769093324, 0.01516077844306046, 3.0, 1.0, 0.657666, 1.972998, 0.1894790521235985, 0.1976066831789769, 50.171543382987, 301.029260297922, 0.004276556927234007, 0.002995871142544512),
  @("MuSCs", "Bgn", "Fgfr3", "MuSCs", 2.0, 1.0, 76.28726950000001, 152.574539, 0.02257007769093324, 0.01516077844306046, 2.0, 1.0, 0.42828, 0.85656, 0.123391035029171, 0.08578923067523865, 32.67231178146, 130.68924712584, 0.002784945246973054, 0.001300631519067899),
  @("MuSCs", "Bgn", "Fgfr3", "Resolving-Mac", 2.0, 1.0, 76.28726950000001, 152.574539, 0.02257007769093324, 0.01516077844306046, 1.0, 0.3333333333333333, 0.05413033333333334, 0.162391, 0.01559539987035126, 0.01626435854882633, 4.129455327124834, 24.776731962749, 0.0003519893866949981, 0.0002465803364772523),
  @("Resolving-Mac", "Bgn", "Fgfr3", "ECs", 3.0, 1.0, 3.227899333333333, 9.683698, 0.0009549947089382317, 0.0009622339405364853, 3.0, 1.0, 2.330840333333333, 6.992521, 0.6715345129768794, 0.7003397275969581, 7.52371795807311, 67.713461622658, 0.0006413119067623321, 0.0006738906557998697),
  @("Resolving-Mac", "Bgn", "Fgfr3", "FAPs", 3.0, 1.0, 3.227899333333333, 9.683698, 0.0009549947089382317, 0.0009622339405364853, 3.0, 1.0, 0.657666, 1.972998, 0.1894790521235985, 0.1976066831789769, 2.122879642956, 19.105916786604, 0.000180951492232668, 0.0001901438574316518),
  @("Resolving-Mac", "Bgn", "Fgfr3", "MuSCs", 3.0, 1.0, 3.227899333333333, 9.683698, 0.0009549947089382317, 0.0009622339405364853, 2.0, 1.0, 0.42828, 0.85656, 0.123391035029171, 0.08578923067523865, 1.38244472648, 8.29466835888, 0.0001178377855832703, 0.0000825493094882284055881),
  @("Resolving-Mac", "Bgn", "Fgfr3", "Resolving-Mac", 3.0, 1.0, 3.227899333333333, 9.683698, 0.0009549947089382317, 0.0009622339405364853, 1.0, 0.3333333333333333, 0.05413033333333334, 0.162391, 0.01559539987035126, 0.01626435854882633, 0.1747272668797778, 1.572545401918, 0.0000148935243599614295724, 0.0000156501178167354301383)
)

$startRow = 2
$r = $startRow
foreach ($row in $rows) {
  $c = 1
  foreach ($val in $row) {
    $ws.Cells.Item($r, $c).Value = $val
    $c = $c + 1
  }
  $r = $r + 1
}

